# Fill in impact type ("tipo_impacto") and impact description ("descricao_impacto")
# for the IOPC source rows (81-150), which previously held "NULL" placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$value = "Derramamento de petróleo"

for ($row = 81; $row -le 150; $row++) {
    $ws.Cells.Item($row, 12).Value = $value   # Column L: tipo_impacto
    $ws.Cells.Item($row, 13).Value = $value   # Column M: descricao_impacto
}
